## Weekly Fruta/Hortaliza update: insert a new price record as row 235
## (pushing the existing rows 235..317 down to 236..318) on the
## "Espinaca" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 235; everything from the old row 235 onward
# (through the old last row 317) shifts down by one, becoming 236..318,
# and the used range grows to A1:R318.
$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with the new weekly record.
$ws.Range("A235").Value = 10
$ws.Range("B235").Value = "Vega Modelo de Temuco"
$ws.Range("C235").Value = "La Araucanía"
$ws.Range("D235").Value = 45229
$ws.Range("E235").Value = 9
$ws.Range("F235").Value = 100112012
$ws.Range("G235").Value = "Espinaca"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 95
$ws.Range("K235").Value = 10000
$ws.Range("L235").Value = 12000
$ws.Range("M235").Value = 11263
$ws.Range("N235").Value = "$/docena de atados"
$ws.Range("O235").Value = "Región de La Araucanía"
$ws.Range("P235").Value = 3754
$ws.Range("Q235").Value = 3
$ws.Range("R235").Value = "Hortaliza"
